$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H39").Value = 300.1
$ws.Range("I39").Value = 385.85715
$ws.Range("K39").Value = 1157.57145
$ws.Range("M39").Value = -861.5714499999999
$ws.Range("H40").Value = 3299.2
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 3624
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 3624
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -3974
$ws.Range("H137").Value = 1536.0588
$ws.Range("J137").Value = 1558.1111
$ws.Range("L137").Value = 4674.3333
$ws.Range("N137").Value = -9774.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 811
$ws.Range("I61").Value = 811
$ws.Range("K61").Value = 811
$ws.Range("M61").Value = -599
$ws.Range("H63").Value = 3750
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814
$ws.Range("H66").Value = 3750
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068
$ws.Range("H132").Value = 1789.8
$ws.Range("I132").Value = 1533
$ws.Range("K132").Value = 4599
$ws.Range("M132").Value = -2069
$ws.Range("H133").Value = 49630
$ws.Range("J133").Value = 49630
$ws.Range("L133").Value = 49630
$ws.Range("N133").Value = -54690
$ws.Range("H136").Value = 811
$ws.Range("I136").Value = 811
$ws.Range("K136").Value = 2433
$ws.Range("M136").Value = 117

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3139.3572
$ws.Range("I86").Value = 2725.375
$ws.Range("J86").Value = 3691.3333
$ws.Range("K86").Value = 2725.375
$ws.Range("L86").Value = 3691.3333
$ws.Range("M86").Value = -1602.375
$ws.Range("N86").Value = -5937.3333
$ws.Range("H89").Value = 3139.3572
$ws.Range("I89").Value = 2725.375
$ws.Range("J89").Value = 3691.3333
$ws.Range("K89").Value = 13626.875
$ws.Range("L89").Value = 18456.6665
$ws.Range("M89").Value = -8010.875
$ws.Range("N89").Value = -29688.6665
$ws.Range("H94").Value = 4979.8
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4979.8
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 4979.8
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -5881.8
$ws.Range("H99").Value = 2194.5386
$ws.Range("I99").Value = 1920.909
$ws.Range("J99").Value = 3699.5
$ws.Range("K99").Value = 1920.909
$ws.Range("L99").Value = 3699.5
$ws.Range("M99").Value = -422.9090000000001
$ws.Range("N99").Value = -6695.5
$ws.Range("H134").Value = 1786.8
$ws.Range("I134").Value = 1786.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5360.4
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2825.4
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H58").Value = 1369.3334
$ws.Range("I58").Value = 1369.3334
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1369.3334
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1166.3334
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value = 2500
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2500
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 12500
$ws.Range("N65").Value = -18740
$ws.Range("H136").Value = 1369.3334
$ws.Range("I136").Value = 1369.3334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4108.0002
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1558.0002
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 28.774193
$ws.Range("I2").Value = 19.307692
$ws.Range("K2").Value = 115.846152
$ws.Range("M2").Value = -2.846151999999989

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1965.5
$ws.Range("I102").Value = 1931
$ws.Range("K102").Value = 1931
$ws.Range("M102").Value = -309
$ws.Range("H132").Value = 2501.3333
$ws.Range("I132").Value = 2699.625
$ws.Range("K132").Value = 8098.875
$ws.Range("M132").Value = -5568.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1202
$ws.Range("J7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("N7").Value = -1724
$ws.Range("H46").Value = 4374.9
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 4062.35
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 4062.35
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -4438.35
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H100").Value = 2500
$ws.Range("I100").Value = 2500
$ws.Range("K100").Value = 2500
$ws.Range("M100").Value = -1959
$ws.Range("H126").Value = 1202
$ws.Range("J126").Value = 1500
$ws.Range("L126").Value = 4500
$ws.Range("N126").Value = -9440
$ws.Range("H136").Value = 4230
$ws.Range("I136").Value = 4230
$ws.Range("K136").Value = 12690
$ws.Range("M136").Value = -10140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 70000
$ws.Range("J46").Value = 70000
$ws.Range("L46").Value = 70000
$ws.Range("N46").Value = -70462
$ws.Range("H113").Value = 429.8
$ws.Range("I113").Value = 429.8
$ws.Range("K113").Value = 1289.4
$ws.Range("M113").Value = 880.5999999999999
$ws.Range("H132").Value = 2599.9092
$ws.Range("I132").Value = 2362.375
$ws.Range("K132").Value = 7087.125
$ws.Range("M132").Value = -4557.125
$ws.Range("H134").Value = 70000
$ws.Range("J134").Value = 70000
$ws.Range("L134").Value = 210000
$ws.Range("N134").Value = -215070
$ws.Range("H136").Value = 896.7222
$ws.Range("I136").Value = 861.2353
$ws.Range("K136").Value = 2583.7059
$ws.Range("M136").Value = -33.70589999999993
